$d = $word.ActiveDocument

function Find-ParagraphByPrefix($doc, $needle) {
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Text.StartsWith($needle)) {
            return $p
        }
    }
    return $null
}

# ------------------------------------------------------------------
# Change 1: merge the "...preliminary assessment of salmonid status
# and trends" (FirstParagraph) paragraph with the following
# ", and 2) recommendations..." (BodyText) paragraph into a single
# FirstParagraph-styled paragraph.
# ------------------------------------------------------------------
$pFirst = Find-ParagraphByPrefix $d "The Juvenile Salmonid and Stream Habitat"
$pSecond = Find-ParagraphByPrefix $d ", and 2) recommendations to guide future monitoring"

# Append the lead sentence of the second paragraph onto the run text of
# the first paragraph (range-scoped Find so only this paragraph is hit).
$pFirst.Range.Find.Execute(
    "Analyses of the database conducted during the summer of 2018 provided 1) a preliminary assessment of salmonid status and trends",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Analyses of the database conducted during the summer of 2018 provided 1) a preliminary assessment of salmonid status and trends, and 2) recommendations to guide future monitoring efforts based on conclusions from the existing database. These analyses and recommendations were presented online at",
    2)

# Remove the now-duplicated lead text from the second paragraph (range-scoped).
$pSecond.Range.Find.Execute(
    ", and 2) recommendations to guide future monitoring efforts based on conclusions from the existing database. These analyses and recommendations were presented online at",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "",
    2)

# Delete the paragraph mark that ends the first paragraph, merging the two
# paragraphs into one (the resulting paragraph inherits the second
# paragraph's mark/formatting, so we restore the FirstParagraph style).
$pFirst = Find-ParagraphByPrefix $d "The Juvenile Salmonid and Stream Habitat"
$markRange = $d.Range($pFirst.Range.End - 1, $pFirst.Range.End)
$markRange.Delete()
$pFirst = Find-ParagraphByPrefix $d "The Juvenile Salmonid and Stream Habitat"
$pFirst.Style = "FirstParagraph"

# ------------------------------------------------------------------
# Change 2: intro sentence before the dataset bullet list.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "Additional habitat datasets for the four watersheds in the existing database are available but were not included in the previous analysis. These datasesets include:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The initial analysis incorporated the existing datasets. Since this time, two additional datasets have been compiled and one dataset is in development. These datasesets include:",
    2)

# ------------------------------------------------------------------
# Change 3: first dataset bullet (flow estimates -> flow values).
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "Continous flow estimates covering the period of record, modelled for all sampling stations with salmonid density",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Estimated mean June and September flow values that will be assigned to each of the core fish monitoring sites. These data are currently in development.",
    2)

# ------------------------------------------------------------------
# Change 4: second dataset bullet (reach scale -> reach segment,
# drop the woody-structure sentence, which moves to its own bullet).
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "Reach scale habitat measurements covering 1/2 mile segments for all locations in the study area. This dataset includes observations of woody structure.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Reach segment habitat measurements covering 1/2 mile segments for all locations in the study area.",
    2)

# ------------------------------------------------------------------
# Change 5: insert a new third dataset bullet for stream wood
# inventories, right after the "Reach segment..." bullet.
# ------------------------------------------------------------------
$pReach = Find-ParagraphByPrefix $d "Reach segment habitat measurements"
$pReach.Range.InsertParagraphAfter()
$pNew = Find-ParagraphByPrefix $d "Reach segment habitat measurements"
$pNew = $pNew.Next()
$pNew.Range.InsertAfter("Stream Wood Inventories for the reach segments. This dataset includes observations of woody structure.")

# ------------------------------------------------------------------
# Change 6: "work completed earlier this year" -> "previous work".
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "This follow-up analysis will build on work completed earlier this year to develop a more comprehensive assessment of habitat trends and associations with salmonid densities. The analysis will be driven by the following questions:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "This follow-up analysis will build on previous work to develop a more comprehensive assessment of habitat trends and associations with salmonid densities. The analysis will be driven by the following questions:",
    2)

# ------------------------------------------------------------------
# Change 7: flow question expanded into a two-part question.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "Are flow changes linked to changes in salmonid density?",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "What habitat variables are associated with changes in juvenile steelhead densities? How critical are changes in flow in describing changes in steelhead densities?",
    2)

# ------------------------------------------------------------------
# Change 8: "woody habitat" -> "stream wood".
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "Is woody habitat associated with salmonid density, and if so, is there variation by location and/or age classes?",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Is stream wood associated with salmonid density, and if so, is there variation by location and/or age classes?",
    2)

# ------------------------------------------------------------------
# Change 9: prepend a sentence about the technical memo.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "Overall, this work will inform future monitoring programs by building on the knowledge gained from the first analysis. This work will also inform conservation activities by identifying habitat conditions that are likely to sustain healthy salmonid populations in the region.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The new analysis will be incorporated into the existing technical memo which will summarize findings, evaluate monitoring methods for their ability to describe trends, and provide recommendations for the monitoring program. Overall, this work will inform future monitoring programs by building on the knowledge gained from the first analysis. This work will also inform conservation activities by identifying habitat conditions that are likely to sustain healthy salmonid populations in the region.",
    2)

# ------------------------------------------------------------------
# Change 10: "will be available" -> "will be provided".
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "A sum of `$5000 will be available to support the continued analysis of habitat data. These funds will be applied directly to data analysis by M. Beck and will also support time between collaborators to discuss progress towards meeting objectives.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "A sum of `$5000 will be provided to support the continued analysis of habitat data. These funds will be applied directly to data analysis by M. Beck and will also support time between collaborators to discuss progress towards meeting objectives.",
    2)
